$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Active Signals"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Active Signals")

# Row 2: USDCHF/BUY -> AUDUSD/BUY (signal style unchanged, stays BUY/green)
$ws1.Range("A2").Value = "2025-07-28 19:44"
$ws1.Range("B2").Value = "AUDUSD"
$ws1.Range("C2").Value = "BUY"
$ws1.Range("D2").Value = 0.6569199999999999
$ws1.Range("E2").Value = 0.65464
$ws1.Range("F2").Value = 0.66323
$ws1.Range("G2").Value = 0.05
$ws1.Range("H2").Value = "'77.0%"
$ws1.Range("I2").Value = 2.77
$ws1.Range("J2").Value = "Active"

# Row 3: USDCHF/SELL -> USDCAD/BUY (signal flips SELL -> BUY, pick up green fill)
$ws1.Range("C2").Copy()
$ws1.Range("C3").PasteSpecial(-4122)
$ws1.Range("A3").Value = "2025-07-28 19:29"
$ws1.Range("B3").Value = "USDCAD"
$ws1.Range("C3").Value = "BUY"
$ws1.Range("D3").Value = 1.36116
$ws1.Range("E3").Value = 1.35649
$ws1.Range("F3").Value = 1.36818
$ws1.Range("G3").Value = 0.09
$ws1.Range("H3").Value = "'68.0%"
$ws1.Range("I3").Value = 1.51
$ws1.Range("J3").Value = "Active"

# Row 4: AUDUSD/SELL -> GBPUSD/SELL (signal stays SELL/red)
$ws1.Range("A4").Value = "2025-07-28 19:24"
$ws1.Range("B4").Value = "GBPUSD"
$ws1.Range("C4").Value = "SELL"
$ws1.Range("D4").Value = 1.27154
$ws1.Range("E4").Value = 1.27653
$ws1.Range("F4").Value = 1.26573
$ws1.Range("G4").Value = 0.08
$ws1.Range("H4").Value = "'91.0%"
$ws1.Range("I4").Value = 1.16
$ws1.Range("J4").Value = "Active"

# Row 5: USDCHF/SELL -> AUDUSD/BUY (signal flips SELL -> BUY, pick up green fill)
$ws1.Range("C2").Copy()
$ws1.Range("C5").PasteSpecial(-4122)
$ws1.Range("A5").Value = "2025-07-28 19:40"
$ws1.Range("B5").Value = "AUDUSD"
$ws1.Range("C5").Value = "BUY"
$ws1.Range("D5").Value = 0.65813
$ws1.Range("E5").Value = 0.65364
$ws1.Range("F5").Value = 0.6622400000000001
$ws1.Range("G5").Value = 0.07000000000000001
$ws1.Range("H5").Value = "'65.0%"
$ws1.Range("I5").Value = 0.92
$ws1.Range("J5").Value = "Active"

# ---------------------------------------------------------------------------
# Sheet "Summary Dashboard"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary Dashboard")
$ws2.Range("B5").Value = 10
$ws2.Range("B6").Value = 5
$ws2.Range("B7").Value = "'79.6%"
$ws2.Range("B8").Value = "'1.99"
$ws2.Range("B9").Value = "2025-07-28 19:25:28"

# ---------------------------------------------------------------------------
# Sheet "Signal History"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Signal History")

# Row 2
$ws3.Range("A2").Value = "2025-07-28 19:44"
$ws3.Range("B2").Value = "AUDUSD"
$ws3.Range("C2").Value = "BUY"
$ws3.Range("D2").Value = 0.6569199999999999
$ws3.Range("E2").Value = 0.65464
$ws3.Range("F2").Value = 0.66323
$ws3.Range("G2").Value = 0.05
$ws3.Range("H2").Value = 0.77
$ws3.Range("I2").Value = 2.77
$ws3.Range("J2").Value = "Active"

# Row 3
$ws3.Range("A3").Value = "2025-07-28 19:29"
$ws3.Range("B3").Value = "USDCAD"
$ws3.Range("C3").Value = "BUY"
$ws3.Range("D3").Value = 1.36116
$ws3.Range("E3").Value = 1.35649
$ws3.Range("F3").Value = 1.36818
$ws3.Range("G3").Value = 0.09
$ws3.Range("H3").Value = 0.68
$ws3.Range("I3").Value = 1.51
$ws3.Range("J3").Value = "Active"

# Row 4
$ws3.Range("A4").Value = "2025-07-28 19:38"
$ws3.Range("B4").Value = "USDCAD"
$ws3.Range("C4").Value = "BUY"
$ws3.Range("D4").Value = 1.36194
$ws3.Range("E4").Value = 1.35947
$ws3.Range("F4").Value = 1.3684
$ws3.Range("G4").Value = 0.1
$ws3.Range("H4").Value = 0.77
$ws3.Range("I4").Value = 2.61
$ws3.Range("J4").Value = "Pending"

# Row 5
$ws3.Range("A5").Value = "2025-07-28 19:24"
$ws3.Range("B5").Value = "GBPUSD"
$ws3.Range("C5").Value = "SELL"
$ws3.Range("D5").Value = 1.27154
$ws3.Range("E5").Value = 1.27653
$ws3.Range("F5").Value = 1.26573
$ws3.Range("G5").Value = 0.08
$ws3.Range("H5").Value = 0.91
$ws3.Range("I5").Value = 1.16
$ws3.Range("J5").Value = "Active"

# Row 6
$ws3.Range("A6").Value = "2025-07-28 19:34"
$ws3.Range("B6").Value = "USDJPY"
$ws3.Range("C6").Value = "BUY"
$ws3.Range("D6").Value = 150.00672
$ws3.Range("E6").Value = 149.69719
$ws3.Range("F6").Value = 150.70286
$ws3.Range("G6").Value = 0.09
$ws3.Range("H6").Value = 0.79
$ws3.Range("I6").Value = 2.25
$ws3.Range("J6").Value = "Filled"

# Row 7
$ws3.Range("A7").Value = "2025-07-28 19:02"
$ws3.Range("B7").Value = "GBPUSD"
$ws3.Range("C7").Value = "BUY"
$ws3.Range("D7").Value = 1.27111
$ws3.Range("E7").Value = 1.26667
$ws3.Range("F7").Value = 1.28048
$ws3.Range("G7").Value = 0.08
$ws3.Range("H7").Value = 0.82
$ws3.Range("I7").Value = 2.11
$ws3.Range("J7").Value = "Pending"

# Row 8
$ws3.Range("A8").Value = "2025-07-28 19:11"
$ws3.Range("B8").Value = "NZDUSD"
$ws3.Range("C8").Value = "BUY"
$ws3.Range("D8").Value = 0.58896
$ws3.Range("E8").Value = 0.5866
$ws3.Range("F8").Value = 0.59733
$ws3.Range("G8").Value = 0.04
$ws3.Range("H8").Value = 0.82
$ws3.Range("I8").Value = 3.56
$ws3.Range("J8").Value = "Filled"

# Row 9
$ws3.Range("A9").Value = "2025-07-28 19:03"
$ws3.Range("B9").Value = "EURUSD"
$ws3.Range("C9").Value = "SELL"
$ws3.Range("D9").Value = 1.10897
$ws3.Range("E9").Value = 1.11237
$ws3.Range("F9").Value = 1.10414
$ws3.Range("G9").Value = 0.09
$ws3.Range("H9").Value = 0.74
$ws3.Range("I9").Value = 1.42
$ws3.Range("J9").Value = "Filled"

# Row 10
$ws3.Range("A10").Value = "2025-07-28 19:02"
$ws3.Range("B10").Value = "AUDUSD"
$ws3.Range("C10").Value = "BUY"
$ws3.Range("D10").Value = 0.65579
$ws3.Range("E10").Value = 0.65228
$ws3.Range("F10").Value = 0.6625799999999999
$ws3.Range("G10").Value = 0.03
$ws3.Range("H10").Value = 0.79
$ws3.Range("I10").Value = 1.94
$ws3.Range("J10").Value = "Filled"

# Row 11
$ws3.Range("A11").Value = "2025-07-28 19:21"
$ws3.Range("B11").Value = "AUDUSD"
$ws3.Range("C11").Value = "SELL"
$ws3.Range("D11").Value = 0.65563
$ws3.Range("E11").Value = 0.6601
$ws3.Range("F11").Value = 0.65027
$ws3.Range("G11").Value = 0.02
$ws3.Range("H11").Value = 0.89
$ws3.Range("I11").Value = 1.2
$ws3.Range("J11").Value = "Filled"

# Row 12
$ws3.Range("A12").Value = "2025-07-28 19:41"
$ws3.Range("B12").Value = "NZDUSD"
$ws3.Range("C12").Value = "BUY"
$ws3.Range("D12").Value = 0.59095
$ws3.Range("E12").Value = 0.59576
$ws3.Range("F12").Value = 0.58178
$ws3.Range("G12").Value = 0.02
$ws3.Range("H12").Value = 0.91
$ws3.Range("I12").Value = 1.91
$ws3.Range("J12").Value = "Pending"

# Row 13
$ws3.Range("A13").Value = "2025-07-28 19:40"
$ws3.Range("B13").Value = "AUDUSD"
$ws3.Range("C13").Value = "BUY"
$ws3.Range("D13").Value = 0.65813
$ws3.Range("E13").Value = 0.65364
$ws3.Range("F13").Value = 0.6622400000000001
$ws3.Range("G13").Value = 0.07000000000000001
$ws3.Range("H13").Value = 0.65
$ws3.Range("I13").Value = 0.92
$ws3.Range("J13").Value = "Active"

# Row 14
$ws3.Range("A14").Value = "2025-07-28 19:24"
$ws3.Range("B14").Value = "GBPUSD"
$ws3.Range("C14").Value = "BUY"
$ws3.Range("D14").Value = 1.27136
$ws3.Range("E14").Value = 1.26819
$ws3.Range("F14").Value = 1.27863
$ws3.Range("G14").Value = 0.06
$ws3.Range("H14").Value = 0.82
$ws3.Range("I14").Value = 2.3
$ws3.Range("J14").Value = "Filled"

# Row 15
$ws3.Range("A15").Value = "2025-07-28 19:24"
$ws3.Range("B15").Value = "USDCHF"
$ws3.Range("C15").Value = "BUY"
$ws3.Range("D15").Value = 0.87943
$ws3.Range("E15").Value = 0.87619
$ws3.Range("F15").Value = 0.88415
$ws3.Range("G15").Value = 0.07000000000000001
$ws3.Range("H15").Value = 0.88
$ws3.Range("I15").Value = 1.45
$ws3.Range("J15").Value = "Filled"

# Row 16
$ws3.Range("A16").Value = "2025-07-28 19:21"
$ws3.Range("B16").Value = "GBPUSD"
$ws3.Range("C16").Value = "SELL"
$ws3.Range("D16").Value = 1.26355
$ws3.Range("E16").Value = 1.26715
$ws3.Range("F16").Value = 1.25387
$ws3.Range("G16").Value = 0.07000000000000001
$ws3.Range("H16").Value = 0.7
$ws3.Range("I16").Value = 2.69
$ws3.Range("J16").Value = "Pending"
